# Apply updated dSF (column F) values after a data repull / recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -4
    "F3"  = 4
    "F4"  = -7
    "F5"  = 2
    "F7"  = -2
    "F8"  = -9
    "F10" = -2
    "F11" = 2
    "F12" = 0
    "F13" = -1
    "F14" = 5
    "F15" = 4
    "F16" = -2
    "F18" = -1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
